# assignment #2 amelia -> eddie
# added two figs, all data already clean
#
# Clean up the "Data" sheet headers: drop the spaces from the
# "Waist Size" and "Baggy Pants" column headers so they read as single
# tokens ("WaistSize" / "BaggyPants").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D1").Value = "WaistSize"
$ws.Range("E1").Value = "BaggyPants"

# Reflect where the user's selection ended up after the edit.
$ws.Activate() | Out-Null
$ws.Range("E1").Select() | Out-Null
